# Auto-generated: update cryptos price/volume columns (D, E) for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.211.22"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "1.830.34"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'237.20"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("D6").Value = "'0.6090"
$ws.Range("E6").Value = "  -3.88%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "'0.07085"
$ws.Range("E8").Value = "  -5.20%  "
$ws.Range("D9").Value = "'0.2815"
$ws.Range("E9").Value = "  -3.09%  "
$ws.Range("D10").Value = "'23.85"
$ws.Range("E10").Value = "  -4.94%  "
$ws.Range("D11").Value = "'0.07645"
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("D12").Value = "1.833.39"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").Value = "'4.808"
$ws.Range("E13").Value = "  -3.62%  "
$ws.Range("D14").Value = "'0.6325"
$ws.Range("E14").Value = "  -6.87%  "
$ws.Range("D15").Value = "'0.000009982"
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("D16").Value = "2.075.24"
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").Value = "'79.45"
$ws.Range("E17").Value = "  -3.14%  "
$ws.Range("D18").Value = "'5.959"
$ws.Range("E18").Value = "  -4.93%  "
$ws.Range("D19").Value = "29.219.20"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").Value = "'229.20"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").Value = "'11.81"
$ws.Range("E21").Value = "  -4.28%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "'7.042"
$ws.Range("E23").Value = "  -5.14%  "
$ws.Range("D24").Value = "'1.004"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").Value = "'155.59"
$ws.Range("D26").Value = "'8.113"
$ws.Range("E26").Value = "  -4.57%  "
$ws.Range("D27").Value = "'0.1300"
$ws.Range("E27").Value = "  -4.24%  "
$ws.Range("D28").Value = "'16.73"
$ws.Range("E28").Value = "  -4.27%  "
$ws.Range("D29").Value = "'0.06733"
$ws.Range("E29").Value = "  +2.90%  "
$ws.Range("D30").Value = "'1.481"
$ws.Range("E30").Value = "  +3.53%  "
$ws.Range("E31").Value = "  -2.01%  "
$ws.Range("D32").Value = "'3.847"
$ws.Range("E32").Value = "  -5.59%  "
$ws.Range("D33").Value = "'3.833"
$ws.Range("E33").Value = "  -5.43%  "
$ws.Range("D34").Value = "'1.132"
$ws.Range("E34").Value = "  -0.79%  "
$ws.Range("D35").Value = "'1.723"
$ws.Range("E35").Value = "  -6.44%  "
$ws.Range("D36").Value = "'0.6546"
$ws.Range("E36").Value = "  -6.25%  "
$ws.Range("D37").Value = "'2.555"
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("D38").Value = "1.237.00"
$ws.Range("E38").Value = "  -0.88%  "
$ws.Range("D39").Value = "'2.761"
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("D40").Value = "'0.01765"
$ws.Range("E40").Value = "  -5.08%  "
$ws.Range("D41").Value = "'6.595"
$ws.Range("E41").Value = "  -2.60%  "
$ws.Range("D42").Value = "'0.9212"
$ws.Range("E42").Value = "  -1.47%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "1.990.03"
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("D45").Value = "'100.84"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").Value = "'63.56"
$ws.Range("E46").Value = "  -2.83%  "
$ws.Range("D47").Value = "'0.00000000117"
$ws.Range("E47").Value = "  -2.44%  "
$ws.Range("D48").Value = "'1.628"
$ws.Range("E48").Value = "  -5.05%  "
$ws.Range("D49").Value = "'8.573"
$ws.Range("E49").Value = "  -4.98%  "
$ws.Range("D50").Value = "'0.1086"
$ws.Range("E50").Value = "  -5.34%  "
$ws.Range("D51").Value = "'6.530"
$ws.Range("E51").Value = "  -7.59%  "
